# ---------------------------------------------------------------------------
# Nueva red de alcantarillado para diseñar
# Se modificaron los datos adaptandolos a la nueva red de prueba
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja 1: INFORMACION_POZOS  (tabla "Table1", columnas No. / Cota / Cabecera)
# ---------------------------------------------------------------------------
$wsPozos = $wb.Worksheets.Item(1)

# Insertar 4 filas nuevas (5 a 8) copiando el formato de la fila 4, para que
# las celdas nuevas conserven el mismo estilo que el resto de la tabla.
$wsPozos.Rows("4").Copy()
$wsPozos.Rows("5:8").Insert(-4121)   # xlShiftDown

# Actualizar los valores existentes con los de la nueva red de prueba.
$wsPozos.Range("B2").Value2 = 1000
$wsPozos.Range("C2").Value2 = $true

$wsPozos.Range("B3").Value2 = 1001
$wsPozos.Range("C3").Value2 = $true

$wsPozos.Range("B4").Value2 = 998
$wsPozos.Range("C4").Value2 = $false

# Filas nuevas.
$wsPozos.Range("A5").Value2 = 4
$wsPozos.Range("B5").Value2 = 997
$wsPozos.Range("C5").Value2 = $false

$wsPozos.Range("A6").Value2 = 5
$wsPozos.Range("B6").Value2 = 997.5
$wsPozos.Range("C6").Value2 = $true

$wsPozos.Range("A7").Value2 = 6
$wsPozos.Range("B7").Value2 = 996
$wsPozos.Range("C7").Value2 = $false

$wsPozos.Range("A8").Value2 = 7
$wsPozos.Range("B8").Value2 = 994
$wsPozos.Range("C8").Value2 = $false

# Extender la tabla para que cubra los nuevos datos.
$loPozos = $wsPozos.ListObjects.Item(1)
$loPozos.Resize($wsPozos.Range("A1:C8"))

# ---------------------------------------------------------------------------
# Hoja 2: INFORMACION_TUBERIAS (tabla "Table3")
# ---------------------------------------------------------------------------
$wsTub = $wb.Worksheets.Item(2)

# Ajustar la altura de la fila de encabezado.
$wsTub.Rows("1").RowHeight = 30

# Cambiar el encabezado de la columna de caudal individual de m3/s a L/s.
$wsTub.Range("D1").Value2 = "Caudal individual" + [char]10 + "(L/s)"

# Insertar 4 filas nuevas (4 a 7) copiando el formato de la fila 3.
$wsTub.Rows("3").Copy()
$wsTub.Rows("4:7").Insert(-4121)   # xlShiftDown

# Actualizar los valores existentes con los de la nueva red de prueba.
$wsTub.Range("B2").Value2 = 1
$wsTub.Range("C2").Value2 = 3
$wsTub.Range("D2").Value2 = 50
$wsTub.Range("E2").Value2 = 70
$wsTub.Range("F2").Value2 = 0.0015

$wsTub.Range("B3").Value2 = 2
$wsTub.Range("C3").Value2 = 3
$wsTub.Range("D3").Value2 = 100
$wsTub.Range("E3").Value2 = 80
$wsTub.Range("F3").Value2 = 0.0015

# Filas nuevas.
$wsTub.Range("A4").Value2 = 3
$wsTub.Range("B4").Value2 = 3
$wsTub.Range("C4").Value2 = 4
$wsTub.Range("D4").Value2 = 150
$wsTub.Range("E4").Value2 = 50
$wsTub.Range("F4").Value2 = 0.0015

$wsTub.Range("A5").Value2 = 4
$wsTub.Range("B5").Value2 = 4
$wsTub.Range("C5").Value2 = 6
$wsTub.Range("D5").Value2 = 100
$wsTub.Range("E5").Value2 = 60
$wsTub.Range("F5").Value2 = 0.0015

$wsTub.Range("A6").Value2 = 5
$wsTub.Range("B6").Value2 = 5
$wsTub.Range("C6").Value2 = 6
$wsTub.Range("D6").Value2 = 150
$wsTub.Range("E6").Value2 = 50
$wsTub.Range("F6").Value2 = 0.0015

$wsTub.Range("A7").Value2 = 6
$wsTub.Range("B7").Value2 = 6
$wsTub.Range("C7").Value2 = 7
$wsTub.Range("D7").Value2 = 50
$wsTub.Range("E7").Value2 = 80
$wsTub.Range("F7").Value2 = 0.0015

# Extender la tabla para que cubra los nuevos datos.
$loTub = $wsTub.ListObjects.Item(1)
$loTub.Resize($wsTub.Range("A1:F7"))

# ---------------------------------------------------------------------------
# Selecciones finales (igual que en el archivo de referencia).
# ---------------------------------------------------------------------------
$wsPozos.Activate()
$wsPozos.Range("B6").Select()

$wsTub.Activate()
$wsTub.Range("C7").Select()
